# Natmi following Dr Hou advice
# Replace the LR-pairs data rows (rows 2-17) on Sheet1 with the updated
# sending/target cluster combinations and recomputed statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs","Cthrc1","Fzd6","ECs",2,1,2.4135925,4.827185,0.1006166576043813,0.07006619804477368,2,1,25.081504,50.163008,0.670801886161467,0.5799583902495139,60.53652994312001,242.14611977248,0.06749384370028151,0.04063547942895059),
    @("ECs","Cthrc1","Fzd6","FAPs",2,1,2.4135925,4.827185,0.1006166576043813,0.07006619804477368,2,0.6666666666666666,0.3140343333333334,0.942103,0.008398811455623189,0.01089210079525609,0.7579509116758334,4.547705470055,0.000845060336514194,0.0007631680914440504),
    @("ECs","Cthrc1","Fzd6","Neutro",2,1,2.4135925,4.827185,0.1006166576043813,0.07006619804477368,3,1,11.39945533333333,34.198366,0.3048770974345635,0.3953835721837834,27.51363989661834,165.08183937971,0.03067571452399108,0.02770302367227904),
    @("ECs","Cthrc1","Fzd6","sCs",2,1,2.4135925,4.827185,0.1006166576043813,0.07006619804477368,2,1,0.5953364999999999,1.190673,0.01592220494834624,0.01376593677144639,1.43689971137625,5.747598845504999,0.00160203904359454,0.0009645268520999951),
    @("FAPs","Cthrc1","Fzd6","ECs",3,1,19.76703633333333,59.301109,0.8240384930780011,0.8607507786564448,2,1,25.081504,50.163008,0.670801886161467,0.5799583902495139,495.7870008626454,2974.722005175872,0.5527665754263761,0.4991996359956073),
    @("FAPs","Cthrc1","Fzd6","FAPs",3,1,19.76703633333333,59.301109,0.8240384930780011,0.8607507786564448,2,0.6666666666666666,0.3140343333333334,0.942103,0.008398811455623189,0.01089210079525609,6.207528076914111,55.867752692227,0.006920943935537986,0.009375384240721164),
    @("FAPs","Cthrc1","Fzd6","Neutro",3,1,19.76703633333333,59.301109,0.8240384930780011,0.8607507786564448,3,1,11.39945533333333,34.198366,0.3048770974345635,0.3953835721837834,225.3334477542105,2028.001029787894,0.2512304639439726,0.3403267176251583),
    @("FAPs","Cthrc1","Fzd6","sCs",3,1,19.76703633333333,59.301109,0.8240384930780011,0.8607507786564448,2,1,0.5953364999999999,1.190673,0.01592220494834624,0.01376593677144639,11.7680382260595,70.60822935635699,0.01312050977211433,0.01184904079495787),
    @("Neutro","Cthrc1","Fzd6","ECs",3,1,1.151594333333333,3.454783,0.04800709843776271,0.05014589452853316,2,1,25.081504,50.163008,0.670801886161467,0.5799583902495139,28.88371787787734,173.302307267264,0.03220325218119044,0.02908253226839),
    @("Neutro","Cthrc1","Fzd6","FAPs",3,1,1.151594333333333,3.454783,0.04800709843776271,0.05014589452853316,2,0.6666666666666666,0.3140343333333334,0.942103,0.008398811455623189,0.01089210079525609,0.3616401587387778,3.254761428649,0.0004032025683103115,0.0005461941376730642),
    @("Neutro","Cthrc1","Fzd6","Neutro",3,1,1.151594333333333,3.454783,0.04800709843776271,0.05014589452853316,3,1,11.39945533333333,34.198366,0.3048770974345635,0.3953835721837834,13.12754816495311,118.147933484578,0.01463626482796046,0.01982686290904268),
    @("Neutro","Cthrc1","Fzd6","sCs",3,1,1.151594333333333,3.454783,0.04800709843776271,0.05014589452853316,2,1,0.5953364999999999,1.190673,0.01592220494834624,0.01376593677144639,0.6855861398265,4.113516838959,0.0007643788603014907,0.0006903052134274069),
    @("sCs","Cthrc1","Fzd6","ECs",2,1,0.655778,1.311556,0.02733775087985481,0.01903712877024833,2,1,25.081504,50.163008,0.670801886161467,0.5799583902495139,16.447898530112,65.79159412044801,0.01833821485361891,0.01104074255656593),
    @("sCs","Cthrc1","Fzd6","FAPs",2,1,0.655778,1.311556,0.02733775087985481,0.01903712877024833,2,0.6666666666666666,0.3140343333333334,0.942103,0.008398811455623189,0.01089210079525609,0.2059368070446667,1.235620842268,0.0002296046152606975,0.0002073543254178145),
    @("sCs","Cthrc1","Fzd6","Neutro",2,1,0.655778,1.311556,0.02733775087985481,0.01903712877024833,3,1,11.39945533333333,34.198366,0.3048770974345635,0.3953835721837834,7.475512019582666,44.853072117496,0.008334654138639319,0.007526967977303462),
    @("sCs","Cthrc1","Fzd6","sCs",2,1,0.655778,1.311556,0.02733775087985481,0.01903712877024833,2,1,0.5953364999999999,1.190673,0.01592220494834624,0.01376593677144639,0.390408579297,1.561634317188,0.0004352772723358812,0.0002620639109611215)
)

# Clear any existing rows below the header before writing the new block,
# since the new table has more rows (16) than the old one (12, excluding header).
$ws.Rows("2:200").ClearContents()

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $rowValues = $data[$i]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowValues[$c]
    }
}

$wb.Save()
